$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 4
    3  = 7
    4  = 2
    5  = 4
    6  = 2
    7  = 2
    8  = 2
    9  = 3
    10 = 3
    11 = 1
    12 = 1
    13 = 0
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
